# Renderer008-Delete template update:
# - C2's old marker ("#! END_ROW") is pushed into the new column D2,
#   and C2 gets a new "delete a never-existing variable" marker.
# - C3's old marker ("#! FINISH") is pushed into the new column D3,
#   and C3 gets a new "#! DELETE" marker.
# - Active selection moves to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing contents of C2/C3 before overwriting them.
$oldC2 = $ws.Range("C2").Value2
$oldC3 = $ws.Range("C3").Value2

$ws.Range("D2").Value = $oldC2
$ws.Range("C2").Value = "#! DELETE NeverExistedVarShouldNotThrowErrors"

$ws.Range("D3").Value = $oldC3
$ws.Range("C3").Value = "#! DELETE"

$ws.Range("B3").Select()
